# Natmi following Dr Hou advice
# Expand the Gdnf-Ret sender/target cluster grid from 3 to 8 result rows
# (FAPs and sCs as sending clusters; ECs, FAPs, Neutro, sCs as targets).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.13541
$ws.Range("H2").Value = 0.40623
$ws.Range("I2").Value = 0.1064658135528677
$ws.Range("J2").Value = 0.1064658135528677
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.594806333333333
$ws.Range("N2").Value = 4.784419
$ws.Range("O2").Value = 0.1541349506635776
$ws.Range("P2").Value = 0.1541349506635776
$ws.Range("Q2").Value = 0.2159527255966667
$ws.Range("R2").Value = 1.94357453037
$ws.Range("S2").Value = 0.01641010291932891
$ws.Range("T2").Value = 0.01641010291932891

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.13541
$ws.Range("H3").Value = 0.40623
$ws.Range("I3").Value = 0.1064658135528677
$ws.Range("J3").Value = 0.1064658135528677
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.655146333333334
$ws.Range("N3").Value = 16.965439
$ws.Range("O3").Value = 0.546558966355358
$ws.Range("P3").Value = 0.5465589663553578
$ws.Range("Q3").Value = 0.7657633649966668
$ws.Range("R3").Value = 6.89187028497
$ws.Range("S3").Value = 0.05818984500763762
$ws.Range("T3").Value = 0.05818984500763761

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.13541
$ws.Range("H4").Value = 0.40623
$ws.Range("I4").Value = 0.1064658135528677
$ws.Range("J4").Value = 0.1064658135528677
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01599833333333333
$ws.Range("N4").Value = 0.047995
$ws.Range("O4").Value = 0.001546208005005081
$ws.Range("P4").Value = 0.001546208005005081
$ws.Range("Q4").Value = 0.002166334316666667
$ws.Range("R4").Value = 0.01949700885
$ws.Range("S4").Value = 0.0001646182931748225
$ws.Range("T4").Value = 0.0001646182931748225

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Ret"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.13541
$ws.Range("H5").Value = 0.40623
$ws.Range("I5").Value = 0.1064658135528677
$ws.Range("J5").Value = 0.1064658135528677
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.080867333333333
$ws.Range("N5").Value = 9.242602
$ws.Range("O5").Value = 0.2977598749760595
$ws.Range("P5").Value = 0.2977598749760594
$ws.Range("Q5").Value = 0.4171802456066667
$ws.Range("R5").Value = 3.75462221046
$ws.Range("S5").Value = 0.03170124733272634
$ws.Range("T5").Value = 0.03170124733272633

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gdnf"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.136453666666667
$ws.Range("H6").Value = 3.409361
$ws.Range("I6").Value = 0.8935341864471323
$ws.Range("J6").Value = 0.8935341864471323
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.594806333333333
$ws.Range("N6").Value = 4.784419
$ws.Range("O6").Value = 0.1541349506635776
$ws.Range("P6").Value = 0.1541349506635776
$ws.Range("Q6").Value = 1.812423505139889
$ws.Range("R6").Value = 16.311811546259
$ws.Range("S6").Value = 0.1377248477442487
$ws.Range("T6").Value = 0.1377248477442487

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Gdnf"
$ws.Range("C7").Value = "Ret"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.136453666666667
$ws.Range("H7").Value = 3.409361
$ws.Range("I7").Value = 0.8935341864471323
$ws.Range("J7").Value = 0.8935341864471323
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.655146333333334
$ws.Range("N7").Value = 16.965439
$ws.Range("O7").Value = 0.546558966355358
$ws.Range("P7").Value = 0.5465589663553578
$ws.Range("Q7").Value = 6.426811786053222
$ws.Range("R7").Value = 57.841306074479
$ws.Range("S7").Value = 0.4883691213477203
$ws.Range("T7").Value = 0.4883691213477203

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gdnf"
$ws.Range("C8").Value = "Ret"
$ws.Range("D8").Value = "Neutro"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.136453666666667
$ws.Range("H8").Value = 3.409361
$ws.Range("I8").Value = 0.8935341864471323
$ws.Range("J8").Value = 0.8935341864471323
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01599833333333333
$ws.Range("N8").Value = 0.047995
$ws.Range("O8").Value = 0.001546208005005081
$ws.Range("P8").Value = 0.001546208005005081
$ws.Range("Q8").Value = 0.01818136457722222
$ws.Range("R8").Value = 0.163632281195
$ws.Range("S8").Value = 0.001381589711830259
$ws.Range("T8").Value = 0.001381589711830259

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gdnf"
$ws.Range("C9").Value = "Ret"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.136453666666667
$ws.Range("H9").Value = 3.409361
$ws.Range("I9").Value = 0.8935341864471323
$ws.Range("J9").Value = 0.8935341864471323
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.080867333333333
$ws.Range("N9").Value = 9.242602
$ws.Range("O9").Value = 0.2977598749760595
$ws.Range("P9").Value = 0.2977598749760594
$ws.Range("Q9").Value = 3.501262977480222
$ws.Range("R9").Value = 31.511366797322
$ws.Range("S9").Value = 0.2660586276433332
$ws.Range("T9").Value = 0.2660586276433331

